# Applies the scheduled-runner market-price refresh to the Gungnir Profits workbook.
# For each affected Leve row (identified by sheet + row), update the price/profit
# columns (H:N) to the latest scraped values. Cells that no longer carry a value
# (e.g. LeveProfitNQ/HQ no longer applicable) are cleared so they round-trip as
# absent <c> elements, matching how untouched blank cells are already stored.
$wb = $excel.ActiveWorkbook

# ALC!row 20
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 8348.333000000001
$ws.Range("I20").Value = 8348.333000000001
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 8348.333000000001
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = -8118.333000000001

# ALC!row 35
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 8348.333000000001
$ws.Range("I35").Value = 8348.333000000001
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 8348.333000000001
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = -7969.333000000001

# ALC!row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2842886.8
$ws.Range("I40").Value = 5683195.5
$ws.Range("J40").Value = 2578.182
$ws.Range("K40").Value = 5683195.5
$ws.Range("L40").Value = 2578.182
$ws.Range("M40").Value = -5683020.5
$ws.Range("N40").Value = -2928.182

# ALC!row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 644.52
$ws.Range("I92").Value = 412.82352
$ws.Range("K92").Value = 412.82352
$ws.Range("M92").Value = 835.1764800000001

# ALC!row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2052.889
$ws.Range("I96").Value = 1994.5
$ws.Range("J96").Value = 2257.25
$ws.Range("K96").Value = 5983.5
$ws.Range("L96").Value = 6771.75
$ws.Range("M96").Value = -4610.5
$ws.Range("N96").Value = -9517.75

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2943.6826
$ws.Range("J138").Value = 4708.1787
$ws.Range("L138").Value = 14124.5361
$ws.Range("N138").Value = -24404.5361

# ARM!row 11
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 30000000
$ws.Range("I11").Value = 30000000
$ws.Range("K11").Value = 30000000
$ws.Range("M11").Value = -29999856

# ARM!row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 9553792
$ws.Range("I45").Value = 11530079
$ws.Range("J45").Value = 1735.3334
$ws.Range("K45").Value = 11530079
$ws.Range("L45").Value = 1735.3334
$ws.Range("M45").Value = -11529702
$ws.Range("N45").Value = -2489.3334

# ARM!row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1556
$ws.Range("I102").Value = 1570
$ws.Range("K102").Value = 1570
$ws.Range("M102").Value = 52

# BSM!row 22
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 229.3
$ws.Range("I22").Value = 161.25
$ws.Range("J22").Value = 501.5
$ws.Range("K22").Value = 161.25
$ws.Range("L22").Value = 501.5
$ws.Range("M22").Value = 11.75
$ws.Range("N22").Value = -847.5

# BSM!row 80
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 4347
$ws.Range("I80").Value = 980.5714
$ws.Range("J80").Value = 6965.3335
$ws.Range("K80").Value = 980.5714
$ws.Range("L80").Value = 6965.3335
$ws.Range("M80").Value = 17.42859999999996
$ws.Range("N80").Value = -8961.333500000001

# BSM!row 83
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 4347
$ws.Range("I83").Value = 980.5714
$ws.Range("J83").Value = 6965.3335
$ws.Range("K83").Value = 4902.857
$ws.Range("L83").Value = 34826.6675
$ws.Range("M83").Value = 89.14300000000003
$ws.Range("N83").Value = -44810.6675

# BSM!row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 532.3889
$ws.Range("I94").Value = 433.30768
$ws.Range("K94").Value = 433.30768
$ws.Range("M94").Value = 17.69232

# BSM!row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1608.4897
$ws.Range("I99").Value = 1224.3939
$ws.Range("J99").Value = 2400.6875
$ws.Range("K99").Value = 1224.3939
$ws.Range("L99").Value = 2400.6875
$ws.Range("M99").Value = 273.6061
$ws.Range("N99").Value = -5396.6875

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 20834170
$ws.Range("I107").Value = 29412228
$ws.Range("J107").Value = 1742.5714
$ws.Range("K107").Value = 29412228
$ws.Range("L107").Value = 1742.5714
$ws.Range("M107").Value = -29410308
$ws.Range("N107").Value = -5582.5714

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 23810298
$ws.Range("I58").Value = 35714908
$ws.Range("J58").Value = 1077.1428
$ws.Range("K58").Value = 35714908
$ws.Range("L58").Value = 1077.1428
$ws.Range("M58").Value = -35714705
$ws.Range("N58").Value = -1483.1428

# CRP!row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 635.6
$ws.Range("I107").Value = 461
$ws.Range("J107").Value = 897.5
$ws.Range("K107").Value = 461
$ws.Range("L107").Value = 897.5
$ws.Range("M107").Value = 1459
$ws.Range("N107").Value = -4737.5

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 23810298
$ws.Range("I136").Value = 35714908
$ws.Range("J136").Value = 1077.1428
$ws.Range("K136").Value = 107144724
$ws.Range("L136").Value = 3231.4284
$ws.Range("M136").Value = -107142174
$ws.Range("N136").Value = -8331.428400000001

# CUL!row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 16132896
$ws.Range("I122").Value = 27778120
$ws.Range("J122").Value = 8738.385
$ws.Range("K122").Value = 250003080
$ws.Range("L122").Value = 78645.465
$ws.Range("M122").Value = -250000630
$ws.Range("N122").Value = -83545.465

# CUL!row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 15152443
$ws.Range("J129").Value = 37038304
$ws.Range("L129").Value = 111114912
$ws.Range("N129").Value = -111124912

# GSM!row 93
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744

# GSM!row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1372.9
$ws.Range("I97").Value = 1271.5
$ws.Range("J97").Value = 1525
$ws.Range("K97").Value = 1271.5
$ws.Range("L97").Value = 1525
$ws.Range("M97").Value = -775.5
$ws.Range("N97").Value = -2517

# GSM!row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 996.4706
$ws.Range("I102").Value = 938.8570999999999
$ws.Range("J102").Value = 1036.8
$ws.Range("K102").Value = 938.8570999999999
$ws.Range("L102").Value = 1036.8
$ws.Range("M102").Value = 683.1429000000001
$ws.Range("N102").Value = -4280.8

# LTW!row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4168054.2
$ws.Range("I46").Value = 6945027
$ws.Range("J46").Value = 2595
$ws.Range("K46").Value = 6945027
$ws.Range("L46").Value = 2595
$ws.Range("M46").Value = -6944839
$ws.Range("N46").Value = -2971

# LTW!row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4106.4
$ws.Range("I100").Value = 4575
$ws.Range("J100").Value = 3936
$ws.Range("K100").Value = 4575
$ws.Range("L100").Value = 3936
$ws.Range("M100").Value = -4034
$ws.Range("N100").Value = -5018

# WVR!row 14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = $null
$ws.Range("N14").Value = 0

# WVR!row 74
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9799.5
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 9799.5
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = $null
$ws.Range("M74").Value = 9799.5
$ws.Range("N74").Value = -11671.5

# WVR!row 77
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H77").Value = 9799.5
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 9799.5
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = $null
$ws.Range("M77").Value = 29398.5
$ws.Range("N77").Value = -38758.5

# WVR!row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 55556460
$ws.Range("I81").Value = 71429170
$ws.Range("K81").Value = 142858340
$ws.Range("M81").Value = -142857279

# WVR!row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 55556460
$ws.Range("I84").Value = 71429170
$ws.Range("K84").Value = 714291700
$ws.Range("M84").Value = -714286396

# WVR!row 86
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 75010910
$ws.Range("J86").Value = 75010910
$ws.Range("L86").Value = 75010910
$ws.Range("N86").Value = -75013156

# WVR!row 89
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H89").Value = 75010910
$ws.Range("J89").Value = 75010910
$ws.Range("L89").Value = 375054550
$ws.Range("N89").Value = -375065782

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15644357
$ws.Range("I132").Value = 18888724
$ws.Range("J132").Value = 12407.909
$ws.Range("K132").Value = 56666172
$ws.Range("L132").Value = 37223.727
$ws.Range("M132").Value = -56663642
$ws.Range("N132").Value = -42283.727
